$p = $ppt.ActivePresentation

# --- Remove the last slide (previously sldId 259 / slide4.xml) ---
# That slide only held empty title/content placeholders, so it is simply
# deleted from the deck.
$p.Slides.Item($p.Slides.Count).Delete()

# --- Slide 3: tidy up the "Sprint Work" bullet text ---
$s3 = $p.Slides.Item(3)
$content = $s3.Shapes.Item(2)
$tr = $content.TextFrame.TextRange

# 1) Merge "nnecessary files " + "were deleted from " into a single run
#    (no text change, just re-typing across the former run boundary).
$full = $tr.Text
$needle = "nnecessary files were deleted from "
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $needle.Length)
    $rng.Text = $needle
}

# 2) Fix the wording: "persistent orthogonality" -> "orthogonal persistence"
$full = $tr.Text
$needle = "persistent orthogonality"
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $needle.Length)
    $rng.Text = "orthogonal persistence"
}

# 3) Merge "Unit tests were created for heap file and UTF-8 hash " +
#    "table file " + "creation." into a single run.
$full = $tr.Text
$needle = "Unit tests were created for heap file and UTF-8 hash table file creation."
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $needle.Length)
    $rng.Text = $needle
}
